$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Mark every test suite row to run ("N" -> "Y") in the Runmode column (C2:C7)
$ws.Range("C3:C7").Value = "Y"
